$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 10 (shifts existing rows 10-13 down to 11-14)
$ws.Rows.Item(10).Insert()

# Populate the newly inserted row 10 with the new match data
$ws.Range("A10").Value = 'ChaXVScT'
$ws.Range("B10").Value = '17/03/2025'
$ws.Range("C10").Value = '15:00'
$ws.Range("D10").Value = 'ROMANIA - SUPERLIGA'
$ws.Range("E10").Value = 'Sepsi Sf. Gheorghe'
$ws.Range("F10").Value = 'Gloria Buzau'
$ws.Range("G10").Value = 1.48
$ws.Range("H10").Value = 4
$ws.Range("I10").Value = 6.25
$ws.Range("J10").Value = 2.05
$ws.Range("K10").Value = 2.2
$ws.Range("L10").Value = 7
$ws.Range("M10").Value = 1.05
$ws.Range("N10").Value = 11
$ws.Range("O10").Value = 1.3
$ws.Range("P10").Value = 3.4
$ws.Range("S10").Value = 1.95
$ws.Range("T10").Value = 1.85
$ws.Range("W10").Value = 3.4
$ws.Range("X10").Value = 1.3
$ws.Range("Y10").Value = 1.4
$ws.Range("Z10").Value = 2.75
$ws.Range("AA10").Value = 2.1
$ws.Range("AB10").Value = 1.67
$ws.Range("AC10").Value = 5
$ws.Range("AD10").Value = 6
$ws.Range("AE10").Value = 8.5
$ws.Range("AF10").Value = 9
$ws.Range("AG10").Value = 13
$ws.Range("AH10").Value = 34
$ws.Range("AI10").Value = 10
$ws.Range("AJ10").Value = 9
$ws.Range("AK10").Value = 23
$ws.Range("AL10").Value = 81
$ws.Range("AM10").Value = 351
$ws.Range("AN10").Value = 15
$ws.Range("AO10").Value = 34
$ws.Range("AP10").Value = 21
$ws.Range("AQ10").Value = 81
$ws.Range("AR10").Value = 51
$ws.Range("AS10").Value = 51

# Apply updated odds values to rows shifted/affected by the edit
# Row 2
$ws.Range("M2").Value = 1.11
$ws.Range("N2").Value = 6.5
$ws.Range("Q2").Value = 2
$ws.Range("R2").Value = 1.85
$ws.Range("U2").Value = 4.4
$ws.Range("W2").Value = 5.5
$ws.Range("X2").Value = 1.14

# Row 3
$ws.Range("H3").Value = 2.88
$ws.Range("J3").Value = 3
$ws.Range("L3").Value = 5
$ws.Range("M3").Value = 1.14
$ws.Range("N3").Value = 5.5
$ws.Range("W3").Value = 7
$ws.Range("X3").Value = 1.1
$ws.Range("AC3").Value = 5
$ws.Range("AO3").Value = 19
$ws.Range("AR3").Value = 41

# Row 4
$ws.Range("G4").Value = 2.5
$ws.Range("H4").Value = 2.9
$ws.Range("I4").Value = 3.2
$ws.Range("J4").Value = 3.4
$ws.Range("K4").Value = 1.91
$ws.Range("M4").Value = 1.13
$ws.Range("N4").Value = 6
$ws.Range("O4").Value = 1.5
$ws.Range("P4").Value = 2.5
$ws.Range("Q4").Value = 1.98
$ws.Range("R4").Value = 1.88
$ws.Range("S4").Value = 2.6
$ws.Range("T4").Value = 1.48
$ws.Range("U4").Value = 4.3
$ws.Range("V4").Value = 1.22
$ws.Range("W4").Value = 5
$ws.Range("X4").Value = 1.17
$ws.Range("Y4").Value = 1.57
$ws.Range("Z4").Value = 2.25
$ws.Range("AA4").Value = 2.05
$ws.Range("AB4").Value = 1.7
$ws.Range("AC4").Value = 6.5
$ws.Range("AF4").Value = 23
$ws.Range("AG4").Value = 23
$ws.Range("AK4").Value = 17
$ws.Range("AL4").Value = 67
$ws.Range("AN4").Value = 7.5
$ws.Range("AO4").Value = 15

# Row 5
$ws.Range("M5").Value = 1.18
$ws.Range("N5").Value = 4.5
$ws.Range("O5").Value = 1.83
$ws.Range("P5").Value = 1.83
$ws.Range("S5").Value = 4
$ws.Range("T5").Value = 1.25
$ws.Range("Y5").Value = 1.85
$ws.Range("Z5").Value = 1.95

# Row 7
$ws.Range("Y7").Value = 1.75
$ws.Range("Z7").Value = 2.05

# Row 8
$ws.Range("M8").Value = 1.07
$ws.Range("N8").Value = 9
$ws.Range("S8").Value = 2.1
$ws.Range("T8").Value = 1.7

# Row 11
$ws.Range("G11").Value = 1.83
$ws.Range("H11").Value = 3.2
$ws.Range("I11").Value = 4.15
$ws.Range("J11").Value = 2.42
$ws.Range("K11").Value = 2.05
$ws.Range("L11").Value = 4.7
$ws.Range("M11").Value = 1.09
$ws.Range("N11").Value = 6.2
$ws.Range("O11").Value = 1.42
$ws.Range("P11").Value = 2.67
$ws.Range("S11").Value = 2.2
$ws.Range("T11").Value = 1.6
$ws.Range("W11").Value = 3.85
$ws.Range("X11").Value = 1.22
$ws.Range("Y11").Value = 1.45
$ws.Range("Z11").Value = 2.55
$ws.Range("AA11").Value = 2.02
$ws.Range("AB11").Value = 1.7
$ws.Range("AC11").Value = 5.8
$ws.Range("AD11").Value = 7.8
$ws.Range("AE11").Value = 8.5
$ws.Range("AF11").Value = 15
$ws.Range("AG11").Value = 16.5
$ws.Range("AH11").Value = 35
$ws.Range("AI11").Value = 6.2
$ws.Range("AJ11").Value = 6.4
$ws.Range("AK11").Value = 18
$ws.Range("AL11").Value = 110
$ws.Range("AN11").Value = 9.5
$ws.Range("AO11").Value = 22
$ws.Range("AP11").Value = 14.5
$ws.Range("AQ11").Value = 70
$ws.Range("AR11").Value = 50
$ws.Range("AS11").Value = 60

# Row 12
$ws.Range("G12").Value = 1.95
$ws.Range("H12").Value = 3.6
$ws.Range("I12").Value = 3.6
$ws.Range("J12").Value = 2.6
$ws.Range("K12").Value = 2.25
$ws.Range("AD12").Value = 10
$ws.Range("AF12").Value = 17
$ws.Range("AP12").Value = 13
$ws.Range("AR12").Value = 29
